$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a number-looking string ("1.", "2.", ...) to be stored as
# genuine text rather than being auto-coerced into a numeric value. We
# round-trip it through a text formula and then "paste values" over itself,
# which keeps the result a plain literal string cell with no style residue.
function Set-TextValue($range, [string]$text) {
    $escaped = $text -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

# --- Move the I/O reference table block from columns P:W (old) to H:O (new) ---
# Row 1 header
$ws.Range("P1").Value = $null
$ws.Range("Q1").Value = $null
$ws.Range("H1").Value = "I/O:"
$ws.Range("I1").Value = "0000|0|000"

# Row 2
$ws.Range("Q2").Value = $null
$ws.Range("I2").Value = "адрес|значение|момент игры"

# Row 8 headers
$ws.Range("P8").Value = $null
$ws.Range("U8").Value = $null
$ws.Range("W8").Value = $null
$ws.Range("H8").Value = "Поле/ адрес:"
$ws.Range("M8").Value = "значениеЖ"
$ws.Range("O8").Value = "момент игры:"

# Row 9
$ws.Range("Q9").Value = $null
$ws.Range("R9").Value = $null
$ws.Range("S9").Value = $null
$ws.Range("U9").Value = $null
$ws.Range("W9").Value = $null
$ws.Range("I9").Value = "0000/0"
$ws.Range("J9").Value = "0100/4"
$ws.Range("K9").Value = "1000/8"
$ws.Range("M9").Value = "0 = нолик"
$ws.Range("O9").Value = "000 - игра идет"

# Row 10
$ws.Range("Q10").Value = $null
$ws.Range("R10").Value = $null
$ws.Range("S10").Value = $null
$ws.Range("U10").Value = $null
$ws.Range("W10").Value = $null
$ws.Range("I10").Value = "0001/1"
$ws.Range("J10").Value = "0101/5"
$ws.Range("K10").Value = "1001/9"
$ws.Range("M10").Value = "1 = крестик"
$ws.Range("O10").Value = "001 - победа игрока(1)"

# Row 11
$ws.Range("Q11").Value = $null
$ws.Range("R11").Value = $null
$ws.Range("S11").Value = $null
$ws.Range("W11").Value = $null
$ws.Range("I11").Value = "0010/2"
$ws.Range("J11").Value = "0110/6"
$ws.Range("K11").Value = "1010/А"
$ws.Range("O11").Value = "010 - победа бота (0)"

# Row 12
$ws.Range("T12").Value = $null
$ws.Range("W12").Value = $null
$ws.Range("L12").Value = "1111/F"
$ws.Range("O12").Value = "011- ничья"

# --- New column A numbering / separators (interleaved with existing rows) ---
Set-TextValue $ws.Range("A2") "1."
$ws.Range("A3").Value = "-"
Set-TextValue $ws.Range("A4") "2."
$ws.Range("A5").Value = "+"
Set-TextValue $ws.Range("A6") "3."
$ws.Range("A7").Value = "-"
Set-TextValue $ws.Range("A8") "4."
$ws.Range("A9").Value = "-"
Set-TextValue $ws.Range("A10") "5."
$ws.Range("A11").Value = "-"
Set-TextValue $ws.Range("A12") "#########################################################################################################################################"

# --- Updated E2 text ---
$ws.Range("E2").Value = "выводить крестих/нолик"

# --- New "Ассемблер" block (rows 13-19) ---
$ws.Range("B13").Value = "Ассемблер : "
$ws.Range("D13").Value = "1. создавать таблицу"
$ws.Range("D14").Value = "2. решать кто делает ход:"
$ws.Range("D15").Value = "           А. человек: считывать из BBB пока не получится корректный адрес"
$ws.Range("D16").Value = "           Б. Мега AI делает корректный ход"
$ws.Range("D17").Value = "3. Отправлять ход в LLL"
$ws.Range("D18").Value = "4. Проверять на конец игры wp/wb/draw"
$ws.Range("D19").Value = "5.если конец то отправить в WIN сигнал иначе на 2 пункт"

# --- Column widths ---
$ws.Columns.Item(4).ColumnWidth = 65.6640625
$ws.Columns.Item(6).ColumnWidth = 70.77734375

# --- Row heights ---
$ws.Rows.Item(4).RowHeight = 17.4
$ws.Rows.Item(5).RowHeight = 15

# --- Selection matches the authored state ---
$ws.Range("D10").Select()
